# Append two new match rows (146 and 147) to the "Azerbaijan Premier League"
# sheet, mirroring the formatting already used on row 145 (bold/bordered,
# centered id in column A; date-formatted cell in column E) but only on the
# cells that actually receive data, so no stray empty <c/> elements get
# written for columns that have no value (as in source row 147 which omits
# FTHG/FTAG/FTR and the last two PL columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 146 --------------------------------------------------------------
$ws.Range("A145").Copy()
$ws.Range("A146").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E145").Copy()
$ws.Range("E146").PasteSpecial(-4122)   # xlPasteFormats

$ws.Cells.Item(146, 1).Value = 144
$ws.Cells.Item(146, 2).Value = 7011638
$ws.Cells.Item(146, 3).Value = "Azerbaijan Premier League"
$ws.Cells.Item(146, 4).Value = "Azerbaijan Premier League"
$ws.Cells.Item(146, 5).Value = 45394.5
$ws.Cells.Item(146, 6).Value = "FK Sumqayit"
$ws.Cells.Item(146, 7).Value = "Zira IK"
$ws.Cells.Item(146, 8).Value = 0
$ws.Cells.Item(146, 9).Value = 0
$ws.Cells.Item(146, 10).Value = "D"
$ws.Cells.Item(146, 11).Value = 2.9
$ws.Cells.Item(146, 12).Value = 3
$ws.Cells.Item(146, 13).Value = 2.3
$ws.Cells.Item(146, 14).Value = 2.875
$ws.Cells.Item(146, 15).Value = 2.8
$ws.Cells.Item(146, 16).Value = 2.45
$ws.Cells.Item(146, 17).Value = 0
$ws.Cells.Item(146, 18).Value = 2.05
$ws.Cells.Item(146, 19).Value = 1.75
$ws.Cells.Item(146, 20).Value = 1.75
$ws.Cells.Item(146, 21).Value = 1.8
$ws.Cells.Item(146, 22).Value = 2
$ws.Cells.Item(146, 23).Value = -1
$ws.Cells.Item(146, 24).Value = 1.8
$ws.Cells.Item(146, 25).Value = -1
$ws.Cells.Item(146, 26).Value = 0
$ws.Cells.Item(146, 27).Value = -0
$ws.Cells.Item(146, 28).Value = -1
$ws.Cells.Item(146, 29).Value = 1

# ---- Row 147 --------------------------------------------------------------
$ws.Range("A145").Copy()
$ws.Range("A147").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E145").Copy()
$ws.Range("E147").PasteSpecial(-4122)   # xlPasteFormats

$ws.Cells.Item(147, 1).Value = 145
$ws.Cells.Item(147, 2).Value = 7011640
$ws.Cells.Item(147, 3).Value = "Azerbaijan Premier League"
$ws.Cells.Item(147, 4).Value = "Azerbaijan Premier League"
$ws.Cells.Item(147, 5).Value = 45396.39583333334
$ws.Cells.Item(147, 6).Value = "Sabail FC"
$ws.Cells.Item(147, 7).Value = "FK Kapaz"
$ws.Cells.Item(147, 11).Value = 1.727
$ws.Cells.Item(147, 12).Value = 3.25
$ws.Cells.Item(147, 13).Value = 4.333
$ws.Cells.Item(147, 14).Value = 1.95
$ws.Cells.Item(147, 15).Value = 3.1
$ws.Cells.Item(147, 16).Value = 3.6
$ws.Cells.Item(147, 17).Value = -0.5
$ws.Cells.Item(147, 18).Value = 1.975
$ws.Cells.Item(147, 19).Value = 1.825
$ws.Cells.Item(147, 20).Value = 2.5
$ws.Cells.Item(147, 21).Value = 1.825
$ws.Cells.Item(147, 22).Value = 1.975
$ws.Cells.Item(147, 23).Value = 0
$ws.Cells.Item(147, 24).Value = 0
$ws.Cells.Item(147, 25).Value = 0
$ws.Cells.Item(147, 26).Value = 0
$ws.Cells.Item(147, 27).Value = 0
